# Insert a new weekly price record as row 33 in the "Arveja Verde" sheet.
# All existing rows from 33 downward shift down by one (old row 33 -> new
# row 34, ..., old row 79 -> new row 80), and the new row 33 holds the
# newly reported data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data down by inserting a fresh row at position 33.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new record's data.
$ws.Range("A33").Value = 5
$ws.Range("B33").Value = "Macroferia Regional de Talca"
$ws.Range("C33").Value = "Maule"
$ws.Range("D33").Value = 44540
$ws.Range("E33").Value = 7
$ws.Range("F33").Value = 100112022
$ws.Range("G33").Value = "Arveja Verde"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 350
$ws.Range("K33").Value = 17000
$ws.Range("L33").Value = 17000
$ws.Range("M33").Value = 17000
$ws.Range("N33").Value = "$/saco 25 kilos"
$ws.Range("O33").Value = "Región de La Araucanía"
$ws.Range("P33").Value = 680
$ws.Range("Q33").Value = 25
$ws.Range("R33").Value = "Hortaliza"
